$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header fix: "unnamed: 1_level_1" -> "total"
$ws.Range("B2").Value = "total"

# 2) Remove the source note row at the bottom (row 41) and the two
#    label-only separator rows ("situação do domicílio" / "grandes
#    regiões e unidades da federação") that had no data of their own.
#    Deleting from the bottom up keeps the remaining row numbers stable
#    while we work.
$ws.Rows("41:41").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
$ws.Rows("8:8").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
$ws.Rows("5:5").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
